$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14: replace the 39 k metal-film resistor Conrad link with the 39 k thick-film SMD resistor link, update price
$ws.Range("C14").Value = "https://www.conrad.de/de/dickschicht-widerstand-39-k-smd-0805-0125-w-1-100-ppmc-royalohm-0805s8f3901t5e-1-st-1208685.html"
$ws.Range("E14").Value = 0.02

# Row 12: replace the 10 k carbon-film resistor Conrad link with the 10 k thick-film SMD resistor link, update price
$ws.Range("C12").Value = "https://www.conrad.de/de/dickschicht-widerstand-10-k-smd-1206-025-w-1-royalohm-1206s4f1002t5e-1-st-1208947.html"
$ws.Range("E12").Value = 0.02

# Row 13: replace the 470 k carbon-film resistor Conrad link with the 470 thick-film SMD resistor link, update price
$ws.Range("C13").Value = "https://www.conrad.de/de/dickschicht-widerstand-470-smd-0805-033-w-1-100-ppmc-royalohm-hp05w3f4700t5e-1-st-1376870.html"
$ws.Range("E13").Value = 0.02

# Row 10: replace the 220 Ohm metal-film resistor Conrad link with the 220 Ohm thick-film SMD resistor link, update price
$ws.Range("C10").Value = "https://www.conrad.de/de/dickschicht-widerstand-220-smd-0805-033-w-1-100-ppmc-royalohm-hp05w3f2200t5e-1-st-1376866.html"
$ws.Range("E10").Value = 0.03

# Row 7: replace the MKP-X2 capacitor Conrad link with the ceramic capacitor Conrad link
$ws.Range("C7").Value = "https://www.conrad.de/de/keramik-kondensator-smd-0603-10-pf-50-v-5-1-st-454117.html"

# Move the active selection from A17 to C10, as saved in the workbook view
$ws.Range("C10").Select()
